$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sort the language rows (rows 2-21) in descending order by 2004 value,
# and drop the Uzbek / Vietnamese rows (previously rows 22-23).
$ws.Range("A2").Value = "English"
$ws.Range("B2").Value = 26.04235644113318
$ws.Range("A3").Value = "Chinese"
$ws.Range("B3").Value = 10.35233516613632
$ws.Range("A4").Value = "Spanish"
$ws.Range("B4").Value = 7.490791254328094
$ws.Range("A5").Value = "Japanese"
$ws.Range("B5").Value = 6.218422985332314
$ws.Range("A6").Value = "German"
$ws.Range("B6").Value = 5.179803403620602
$ws.Range("A7").Value = "Arabic"
$ws.Range("B7").Value = 5.076453270470317
$ws.Range("A8").Value = "Portuguese"
$ws.Range("B8").Value = 3.561159143165131
$ws.Range("A9").Value = "Russian"
$ws.Range("B9").Value = 3.4739989100432
$ws.Range("A10").Value = "French"
$ws.Range("B10").Value = 3.269294161405805
$ws.Range("A11").Value = "Italian"
$ws.Range("B11").Value = 2.963023247149283
$ws.Range("A12").Value = "Malay-Indonesian"
$ws.Range("B12").Value = 2.644103529477143
$ws.Range("A13").Value = "Korean"
$ws.Range("B13").Value = 1.660300691674208
$ws.Range("A14").Value = "Persian"
$ws.Range("B14").Value = 1.56994984051646
$ws.Range("A15").Value = "Dutch"
$ws.Range("B15").Value = 1.515150462663296
$ws.Range("A16").Value = "Turkish"
$ws.Range("B16").Value = 1.400670266597371
$ws.Range("A17").Value = "Thai"
$ws.Range("B17").Value = 0.986213562068493
$ws.Range("A18").Value = "Polish"
$ws.Range("B18").Value = 0.8636613076071532
$ws.Range("A19").Value = "Urdu"
$ws.Range("B19").Value = 0.8250723533026133
$ws.Range("A20").Value = "Swedish"
$ws.Range("B20").Value = 0.496055063657645
$ws.Range("A21").Value = "Bengali"
$ws.Range("B21").Value = 0.4519010778952846

# Remove the now-obsolete trailing rows (Uzbek, Vietnamese) and shrink the
# used range / dimension accordingly.
$ws.Range("A22:B23").Delete() | Out-Null
